# "Generate Report for Archive"
#
# 1. Every cell whose status text reads "Ready for handoff" moves to
#    "In Translation" (Overview!E2:E3/F2:F3, zh-cn!C2:C3, de-de!C2:C3 —
#    they all shared one sst entry, so touch them all so no stale
#    reference to the old text is left behind).
# 2. The "status" column narrows on the three sheets (Overview columns
#    E & F, zh-cn column C, de-de column C) from ~17.22 chars to
#    ~13.41 chars.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1. status text -------------------------------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- 2. column widths -------------------------------------------------
# ColumnWidth is quantised to the sheet's pixel grid, so feed it the
# character width whose rounded pixel width lands closest to the
# 13.4101845877511 stored-width target used across the sheets.
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

$zhcn.Range("C1").ColumnWidth = 12.5

$dede.Range("C1").ColumnWidth = 12.5
